# Auto-generated script applying scheduled market-price refresh updates
# to the Hyperion_Profits leve-profit tracking workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2741.5
$ws.Range("I32").Value = 2489.8
$ws.Range("J32").Value = 4000
$ws.Range("K32").Value = 2489.8
$ws.Range("L32").Value = 4000
$ws.Range("M32").Value = -2163.8
$ws.Range("N32").Value = -4652
$ws.Range("H33").Value = 512.3077
$ws.Range("I33").Value = 521.0625
$ws.Range("K33").Value = 521.0625
$ws.Range("M33").Value = -292.0625
$ws.Range("H64").Value = 7374.9165
$ws.Range("J64").Value = 7454.4546
$ws.Range("L64").Value = 7454.4546
$ws.Range("N64").Value = -7950.4546
$ws.Range("H67").Value = 7374.9165
$ws.Range("J67").Value = 7454.4546
$ws.Range("L67").Value = 7454.4546
$ws.Range("N67").Value = -9170.454600000001
$ws.Range("H70").Value = 6522.4
$ws.Range("J70").Value = 6349.857
$ws.Range("L70").Value = 19049.571
$ws.Range("N70").Value = -19589.571
$ws.Range("H73").Value = 6522.4
$ws.Range("J73").Value = 6349.857
$ws.Range("L73").Value = 19049.571
$ws.Range("N73").Value = -20921.571
$ws.Range("H74").Value = 4998.875
$ws.Range("I74").Value = 2416.5
$ws.Range("K74").Value = 2416.5
$ws.Range("M74").Value = -1480.5
$ws.Range("H76").Value = 4030.5833
$ws.Range("I76").Value = 3896.8572
$ws.Range("K76").Value = 3896.8572
$ws.Range("M76").Value = -3581.8572
$ws.Range("H77").Value = 4998.875
$ws.Range("I77").Value = 2416.5
$ws.Range("K77").Value = 12082.5
$ws.Range("M77").Value = -7402.5
$ws.Range("H79").Value = 4030.5833
$ws.Range("I79").Value = 3896.8572
$ws.Range("K79").Value = 3896.8572
$ws.Range("M79").Value = -2804.8572
$ws.Range("H80").Value = 4134.467
$ws.Range("J80").Value = 6164.1113
$ws.Range("L80").Value = 18492.3339
$ws.Range("N80").Value = -20488.3339
$ws.Range("H83").Value = 4134.467
$ws.Range("J83").Value = 6164.1113
$ws.Range("L83").Value = 55477.00169999999
$ws.Range("N83").Value = -65461.00169999999
$ws.Range("H88").Value = 2965.3
$ws.Range("J88").Value = 3116.6072
$ws.Range("L88").Value = 3116.6072
$ws.Range("N88").Value = -3928.6072
$ws.Range("H91").Value = 2965.3
$ws.Range("J91").Value = 3116.6072
$ws.Range("L91").Value = 3116.6072
$ws.Range("N91").Value = -5924.6072
$ws.Range("H92").Value = 1177.7812
$ws.Range("I92").Value = 353.84616
$ws.Range("J92").Value = 4748.1665
$ws.Range("K92").Value = 353.84616
$ws.Range("L92").Value = 4748.1665
$ws.Range("M92").Value = 894.1538399999999
$ws.Range("N92").Value = -7244.1665
$ws.Range("H98").Value = 1292.5807
$ws.Range("I98").Value = 1305.963
$ws.Range("K98").Value = 1305.963
$ws.Range("M98").Value = 192.037
$ws.Range("H100").Value = 3006.2
$ws.Range("I100").Value = 3678.25
$ws.Range("K100").Value = 3678.25
$ws.Range("M100").Value = -3137.25
$ws.Range("H122").Value = 1292.5807
$ws.Range("I122").Value = 1305.963
$ws.Range("K122").Value = 3917.889
$ws.Range("M122").Value = -1467.889
$ws.Range("H133").Value = 124197.5
$ws.Range("J133").Value = 124197.5
$ws.Range("L133").Value = 124197.5
$ws.Range("N133").Value = -134317.5
$ws.Range("H136").Value = 85000
$ws.Range("J136").Value = 85000
$ws.Range("L136").Value = 85000
$ws.Range("N136").Value = -95200
$ws.Range("H137").Value = 63804.207
$ws.Range("J137").Value = 3868
$ws.Range("L137").Value = 11604
$ws.Range("N137").Value = -16704
$ws.Range("H138").Value = 2828.2917
$ws.Range("I138").Value = 1277.2609
$ws.Range("K138").Value = 3831.7827
$ws.Range("M138").Value = 1308.2173
$ws.Range("H141").Value = 6461.72
$ws.Range("I141").Value = 6461.72
$ws.Range("K141").Value = 19385.16
$ws.Range("M141").Value = -14205.16

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 5000
$ws.Range("I6").Value = 5000
$ws.Range("K6").Value = 5000
$ws.Range("M6").Value = -4827
$ws.Range("H32").Value = 11165.1455
$ws.Range("I32").Value = 6493.641
$ws.Range("J32").Value = 22551.938
$ws.Range("K32").Value = 6493.641
$ws.Range("L32").Value = 22551.938
$ws.Range("M32").Value = -6206.641
$ws.Range("N32").Value = -23125.938
$ws.Range("H45").Value = 6157430.5
$ws.Range("I45").Value = 7694538.5
$ws.Range("K45").Value = 7694538.5
$ws.Range("M45").Value = -7694161.5
$ws.Range("H61").Value = 2375.5833
$ws.Range("I61").Value = 2261.6956
$ws.Range("K61").Value = 2261.6956
$ws.Range("M61").Value = -2049.6956
$ws.Range("H63").Value = 5810.6
$ws.Range("I63").Value = 2993.3333
$ws.Range("J63").Value = 7688.778
$ws.Range("K63").Value = 2993.3333
$ws.Range("L63").Value = 7688.778
$ws.Range("M63").Value = -2307.3333
$ws.Range("N63").Value = -9060.778
$ws.Range("H66").Value = 5810.6
$ws.Range("I66").Value = 2993.3333
$ws.Range("J66").Value = 7688.778
$ws.Range("K66").Value = 14966.6665
$ws.Range("L66").Value = 38443.89
$ws.Range("M66").Value = -11534.6665
$ws.Range("N66").Value = -45307.89
$ws.Range("H74").Value = 19425.295
$ws.Range("I74").Value = 1257.7188
$ws.Range("K74").Value = 1257.7188
$ws.Range("M74").Value = -383.7188000000001
$ws.Range("H77").Value = 19425.295
$ws.Range("I77").Value = 1257.7188
$ws.Range("K77").Value = 6288.594000000001
$ws.Range("M77").Value = -1920.594000000001
$ws.Range("H88").Value = 1456.3125
$ws.Range("J88").Value = 1311.7778
$ws.Range("L88").Value = 1311.7778
$ws.Range("N88").Value = -2123.7778
$ws.Range("H91").Value = 1456.3125
$ws.Range("J91").Value = 1311.7778
$ws.Range("L91").Value = 1311.7778
$ws.Range("N91").Value = -4119.7778
$ws.Range("H122").Value = 567268
$ws.Range("I122").Value = 3144.5908
$ws.Range("K122").Value = 9433.7724
$ws.Range("M122").Value = -6983.7724
$ws.Range("H132").Value = 2483.375
$ws.Range("I132").Value = 2063.6667
$ws.Range("J132").Value = 5421.3335
$ws.Range("K132").Value = 6191.000100000001
$ws.Range("L132").Value = 16264.0005
$ws.Range("M132").Value = -3661.000100000001
$ws.Range("N132").Value = -21324.0005
$ws.Range("H136").Value = 2375.5833
$ws.Range("I136").Value = 2261.6956
$ws.Range("K136").Value = 6785.0868
$ws.Range("M136").Value = -4235.0868

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 15154062
$ws.Range("I20").Value = 33335582
$ws.Range("J20").Value = 2796.5
$ws.Range("K20").Value = 33335582
$ws.Range("L20").Value = 2796.5
$ws.Range("M20").Value = -33335335
$ws.Range("N20").Value = -3290.5
$ws.Range("H86").Value = 6251300.5
$ws.Range("I86").Value = 12501162
$ws.Range("J86").Value = 1438.75
$ws.Range("K86").Value = 12501162
$ws.Range("L86").Value = 1438.75
$ws.Range("M86").Value = -12500039
$ws.Range("N86").Value = -3684.75
$ws.Range("H89").Value = 6251300.5
$ws.Range("I89").Value = 12501162
$ws.Range("J89").Value = 1438.75
$ws.Range("K89").Value = 62505810
$ws.Range("L89").Value = 7193.75
$ws.Range("M89").Value = -62500194
$ws.Range("N89").Value = -18425.75
$ws.Range("H96").Value = 18242.4
$ws.Range("I96").Value = 18242.4
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 18242.4
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -15496.4
$ws.Range("N96").ClearContents()
$ws.Range("H99").Value = 5956059
$ws.Range("I99").Value = 10207734
$ws.Range("J99").Value = 3713.9
$ws.Range("K99").Value = 10207734
$ws.Range("L99").Value = 3713.9
$ws.Range("M99").Value = -10206236
$ws.Range("N99").Value = -6709.9
$ws.Range("H105").Value = 6250814.5
$ws.Range("I105").Value = 7813255.5
$ws.Range("K105").Value = 7813255.5
$ws.Range("M105").Value = -7811508.5
$ws.Range("H134").Value = 2578.3833
$ws.Range("I134").Value = 1151.3954
$ws.Range("K134").Value = 3454.1862
$ws.Range("M134").Value = -919.1862000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 702.3570999999999
$ws.Range("I22").Value = 589
$ws.Range("J22").Value = 906.4
$ws.Range("K22").Value = 589
$ws.Range("L22").Value = 906.4
$ws.Range("M22").Value = -239
$ws.Range("N22").Value = -1606.4
$ws.Range("H31").Value = 20704.482
$ws.Range("I31").Value = 2405.9167
$ws.Range("J31").Value = 25695
$ws.Range("K31").Value = 2405.9167
$ws.Range("L31").Value = 25695
$ws.Range("M31").Value = -2110.9167
$ws.Range("N31").Value = -26285
$ws.Range("H34").Value = 20704.482
$ws.Range("I34").Value = 2405.9167
$ws.Range("J34").Value = 25695
$ws.Range("K34").Value = 2405.9167
$ws.Range("L34").Value = 25695
$ws.Range("M34").Value = -2203.9167
$ws.Range("N34").Value = -26099
$ws.Range("H99").Value = 3204.158
$ws.Range("I99").Value = 2725.6
$ws.Range("K99").Value = 2725.6
$ws.Range("M99").Value = -1227.6
$ws.Range("H105").Value = 1104.5
$ws.Range("I105").Value = 1104.5
$ws.Range("K105").Value = 1104.5
$ws.Range("M105").Value = 642.5
$ws.Range("H126").Value = 3204.158
$ws.Range("I126").Value = 2725.6
$ws.Range("K126").Value = 8176.799999999999
$ws.Range("M126").Value = -5706.799999999999
$ws.Range("H132").Value = 50656.273
$ws.Range("I132").Value = 33253.516
$ws.Range("K132").Value = 99760.54800000001
$ws.Range("M132").Value = -97230.54800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 20217960
$ws.Range("I4").Value = 23095142
$ws.Range("K4").Value = 69285426
$ws.Range("M4").Value = -69285314
$ws.Range("H7").Value = 4001.5
$ws.Range("I7").Value = 3668
$ws.Range("J7").Value = 5002
$ws.Range("K7").Value = 11004
$ws.Range("L7").Value = 15006
$ws.Range("M7").Value = -10892
$ws.Range("N7").Value = -15230
$ws.Range("H11").Value = 8478.615
$ws.Range("J11").Value = 8161.8184
$ws.Range("L11").Value = 24485.4552
$ws.Range("N11").Value = -24765.4552
$ws.Range("H31").Value = 1002
$ws.Range("J31").Value = 1002
$ws.Range("L31").Value = 3006
$ws.Range("N31").Value = -3582
$ws.Range("H37").Value = 50725.25
$ws.Range("J37").Value = 50725.25
$ws.Range("L37").Value = 152175.75
$ws.Range("N37").Value = -152399.75
$ws.Range("H40").Value = 22.5
$ws.Range("I40").Value = 22.5
$ws.Range("K40").Value = 90
$ws.Range("M40").Value = -21
$ws.Range("H56").Value = 125003000
$ws.Range("I56").Value = 125003000
$ws.Range("K56").Value = 125003000
$ws.Range("M56").Value = -125002470
$ws.Range("H68").Value = 1609.4
$ws.Range("J68").Value = 1609.4
$ws.Range("L68").Value = 4828.200000000001
$ws.Range("N68").Value = -6450.200000000001
$ws.Range("H71").Value = 1609.4
$ws.Range("J71").Value = 1609.4
$ws.Range("L71").Value = 14484.6
$ws.Range("N71").Value = -22596.6
$ws.Range("H86").Value = 244.2
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 244.2
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H113").Value = 3114.5293
$ws.Range("I113").Value = 12598.75
$ws.Range("J113").Value = 1849.9667
$ws.Range("K113").Value = 37796.25
$ws.Range("L113").Value = 5549.9001
$ws.Range("M113").Value = -35626.25
$ws.Range("N113").Value = -9889.900099999999
$ws.Range("H129").Value = 1665.4
$ws.Range("J129").Value = 1704.0834
$ws.Range("L129").Value = 5112.2502
$ws.Range("N129").Value = -15112.2502
$ws.Range("H137").Value = 2308.5386
$ws.Range("I137").Value = 1316
$ws.Range("K137").Value = 3948
$ws.Range("M137").Value = 1152
$ws.Range("H138").Value = 2890.5454
$ws.Range("I138").Value = 2809.6
$ws.Range("K138").Value = 8428.799999999999
$ws.Range("M138").Value = -3288.799999999999
$ws.Range("H140").Value = 3439.087
$ws.Range("I140").Value = 2676.7646
$ws.Range("J140").Value = 5599
$ws.Range("K140").Value = 8030.293799999999
$ws.Range("L140").Value = 16797
$ws.Range("M140").Value = -2850.293799999999
$ws.Range("N140").Value = -27157

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 344
$ws.Range("J3").Value = 344
$ws.Range("L3").Value = 344
$ws.Range("N3").Value = -576
$ws.Range("H42").Value = 36296
$ws.Range("J42").Value = 36296
$ws.Range("L42").Value = 36296
$ws.Range("N42").Value = -37266
$ws.Range("H45").Value = 28637.691
$ws.Range("J45").Value = 41372.855
$ws.Range("L45").Value = 41372.855
$ws.Range("N45").Value = -42490.855
$ws.Range("H70").Value = 25006362
$ws.Range("I70").Value = 28577414
$ws.Range("K70").Value = 28577414
$ws.Range("M70").Value = -28577144
$ws.Range("H73").Value = 25006362
$ws.Range("I73").Value = 28577414
$ws.Range("K73").Value = 28577414
$ws.Range("M73").Value = -28576478
$ws.Range("H80").Value = 34968600
$ws.Range("I80").Value = 47683076
$ws.Range("K80").Value = 47683076
$ws.Range("M80").Value = -47682078
$ws.Range("H83").Value = 34968600
$ws.Range("I83").Value = 47683076
$ws.Range("K83").Value = 238415380
$ws.Range("M83").Value = -238410388
$ws.Range("H97").Value = 2382037
$ws.Range("I97").Value = 4762661
$ws.Range("J97").Value = 1413.2
$ws.Range("K97").Value = 4762661
$ws.Range("L97").Value = 1413.2
$ws.Range("M97").Value = -4762165
$ws.Range("N97").Value = -2405.2
$ws.Range("H115").Value = 36296
$ws.Range("J115").Value = 36296
$ws.Range("L115").Value = 36296
$ws.Range("N115").Value = -38646
$ws.Range("H132").Value = 3776.25
$ws.Range("I132").Value = 3579.0908
$ws.Range("K132").Value = 10737.2724
$ws.Range("M132").Value = -8207.2724
$ws.Range("H134").Value = 48565.332
$ws.Range("J134").Value = 48565.332
$ws.Range("L134").Value = 145695.996
$ws.Range("N134").Value = -150765.996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7924.6665
$ws.Range("J46").Value = 7924.6665
$ws.Range("L46").Value = 7924.6665
$ws.Range("N46").Value = -8300.666499999999
$ws.Range("H68").Value = 1000
$ws.Range("J68").Value = 1000
$ws.Range("L68").Value = 1000
$ws.Range("N68").Value = -2498
$ws.Range("H71").Value = 1000
$ws.Range("J71").Value = 1000
$ws.Range("L71").Value = 5000
$ws.Range("N71").Value = -12488
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H132").Value = 6299.7026
$ws.Range("I132").Value = 6423.0938
$ws.Range("K132").Value = 19269.2814
$ws.Range("M132").Value = -16739.2814

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 7668.6665
$ws.Range("I14").Value = 6000
$ws.Range("J14").Value = 8503
$ws.Range("K14").Value = 6000
$ws.Range("L14").Value = 8503
$ws.Range("M14").Value = -5832
$ws.Range("N14").Value = -8839
$ws.Range("H107").Value = 71430056
$ws.Range("I107").Value = 76923520
$ws.Range("K107").Value = 230770560
$ws.Range("M107").Value = -230768640
$ws.Range("H109").Value = 67997.5
$ws.Range("J109").Value = 67997.5
$ws.Range("L109").Value = 67997.5
$ws.Range("N109").Value = -70771.5
$ws.Range("H126").Value = 2066.6191
$ws.Range("I126").Value = 2146.3635
$ws.Range("K126").Value = 6439.0905
$ws.Range("M126").Value = -3969.0905
$ws.Range("H132").Value = 17261424
$ws.Range("I132").Value = 19612824
$ws.Range("K132").Value = 58838472
$ws.Range("M132").Value = -58835942
$ws.Range("H136").Value = 1910.5106
$ws.Range("I136").Value = 1478
$ws.Range("J136").Value = 3041.6924
$ws.Range("K136").Value = 4434
$ws.Range("L136").Value = 9125.0772
$ws.Range("M136").Value = -1884
$ws.Range("N136").Value = -14225.0772
$ws.Range("H139").Value = 160976.67
$ws.Range("J139").Value = 160976.67
$ws.Range("L139").Value = 160976.67
$ws.Range("N139").Value = -171256.67
$ws.Range("H141").Value = 63333
$ws.Range("J141").Value = 63333
$ws.Range("L141").Value = 63333
$ws.Range("N141").Value = -73693
